# Inclusao de lib nltk na wordcloud
# Updates the "Planilha1" keyword list: normalizes existing keywords to
# lowercase / corrected spelling, drops a couple of unrelated entries
# (Jesus / jesus abençoado / mensagem / família) and appends the new
# NLTK-driven keyword list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "termo",
    "gerir",
    "transferência",
    "depósito",
    "contrato",
    "dossiê",
    "caixinha",
    "acordo",
    "dinheiro",
    "porã",
    "propina",
    "café",
    "cafezinho",
    "prefeito",
    "governador",
    "deputado",
    "patrão"
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

# Move selection / view to just past the last populated row, matching the
# author's saved cursor position, and bump the zoom level used on save.
$ws.Range("A18").Select()
$excel.ActiveWindow.Zoom = 205
